$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = [double]"22.75000000000012"
$ws.Range("H2").Value = [double]"6.661338147750939e-16"
$ws.Range("I2").Value = [double]"6.661338147750939e-16"
$ws.Range("L2").Value = [double]"47.25970911196823"
$ws.Range("M2").Value = "[39.2551837557467, 55.264234468189755]"
$ws.Range("N2").Value = [double]"1.77635683940025e-15"
$ws.Range("O2").Value = [double]"1.77635683940025e-15"
$ws.Range("P2").Value = [double]"1.641552918091964"
$ws.Range("Q2").Value = "[1.4654476241970409, 1.8176582119868874]"
$ws.Range("T2").Value = [double]"57.10540656633878"
$ws.Range("U2").Value = "[52.01344501731727, 62.19736811536029]"
$ws.Range("X2").Value = [double]"16.80630630630639"
$ws.Range("Y2").Value = [double]"16.16866866866875"
$ws.Range("Z2").Value = [double]"17.44394394394403"
$ws.Range("F3").Value = [double]"22.75000000000012"
$ws.Range("H3").Value = [double]"1.110223024625157e-16"
$ws.Range("I3").Value = [double]"1.110223024625157e-16"
$ws.Range("L3").Value = [double]"55.40776038738573"
$ws.Range("M3").Value = "[48.27588598998206, 62.5396347847894]"
$ws.Range("N3").Value = [double]"0"
$ws.Range("O3").Value = [double]"0"
$ws.Range("P3").Value = [double]"1.893131909370426"
$ws.Range("Q3").Value = "[1.7547634641672714, 2.0315003545735797]"
$ws.Range("T3").Value = [double]"52.2649300560555"
$ws.Range("U3").Value = "[48.26483082077008, 56.265029291340916]"
$ws.Range("X3").Value = [double]"15.89539539539548"
$ws.Range("Y3").Value = [double]"15.39439439439447"
$ws.Range("Z3").Value = [double]"16.39639639639648"
$ws.Range("F4").Value = [double]"22.75000000000012"
$ws.Range("H4").Value = [double]"7.223333042816193e-12"
$ws.Range("I4").Value = [double]"7.223333042816193e-12"
$ws.Range("L4").Value = [double]"43.22628322416569"
$ws.Range("M4").Value = "[31.080010509780216, 55.372555938551166]"
$ws.Range("N4").Value = [double]"5.748925335780086e-09"
$ws.Range("O4").Value = [double]"5.748925335780086e-09"
$ws.Range("P4").Value = [double]"2.245342497160272"
$ws.Range("Q4").Value = "[1.9560266571900398, 2.5346583371305043]"
$ws.Range("T4").Value = [double]"54.622581185936"
$ws.Range("U4").Value = "[48.32984294286025, 60.915319429011745]"
$ws.Range("X4").Value = [double]"14.6201201201202"
$ws.Range("Y4").Value = [double]"13.57257257257264"
$ws.Range("Z4").Value = [double]"15.66766766766776"
$ws.Range("F5").Value = [double]"22.75000000000012"
$ws.Range("L5").Value = [double]"50.5175835339561"
$ws.Range("M5").Value = "[44.12703923614342, 56.90812783176878]"
$ws.Range("N5").Value = [double]"0"
$ws.Range("O5").Value = [double]"0"
$ws.Range("P5").Value = [double]"2.597553084950119"
$ws.Range("Q5").Value = "[2.4591846397469648, 2.735921530153274]"
$ws.Range("T5").Value = [double]"52.73882238508907"
$ws.Range("U5").Value = "[49.28742706748364, 56.190217702694504]"
$ws.Range("X5").Value = [double]"13.34484484484491"
$ws.Range("Y5").Value = [double]"12.84384384384391"
$ws.Range("Z5").Value = [double]"13.84584584584592"
$ws.Range("F6").Value = [double]"22.75000000000012"
$ws.Range("L6").Value = [double]"51.03704081005723"
$ws.Range("M6").Value = "[42.93991255747383, 59.13416906264064]"
$ws.Range("N6").Value = [double]"2.220446049250313e-16"
$ws.Range("O6").Value = [double]"2.220446049250313e-16"
$ws.Range("P6").Value = [double]"3.050395269251351"
$ws.Range("Q6").Value = "[2.874289975356427, 3.226500563146275]"
$ws.Range("T6").Value = [double]"51.0844885572775"
$ws.Range("U6").Value = "[45.82449557597006, 56.344481538584944]"
$ws.Range("X6").Value = [double]"11.70520520520526"
$ws.Range("Y6").Value = [double]"11.06756756756762"
$ws.Range("Z6").Value = [double]"12.34284284284291"
$ws.Range("F7").Value = [double]"22.75000000000012"
$ws.Range("H7").Value = [double]"1.110223024625157e-16"
$ws.Range("I7").Value = [double]"1.110223024625157e-16"
$ws.Range("L7").Value = [double]"50.36766005569682"
$ws.Range("M7").Value = "[41.73989426578246, 58.99542584561117]"
$ws.Range("N7").Value = [double]"2.442490654175344e-15"
$ws.Range("O7").Value = [double]"2.442490654175344e-15"
$ws.Range("P7").Value = [double]"-3.031526844905466"
$ws.Range("Q7").Value = "[-3.2202110883643122, -2.8428426014466197]"
$ws.Range("T7").Value = [double]"50.9675810963374"
$ws.Range("U7").Value = "[45.540874056514, 56.394288136160796]"
$ws.Range("X7").Value = [double]"10.97647647647653"
$ws.Range("Y7").Value = [double]"10.29329329329335"
$ws.Range("Z7").Value = [double]"11.65965965965972"
$ws.Range("F8").Value = [double]"25.8300000000006"
$ws.Range("H8").Value = [double]"5.184741524999481e-13"
$ws.Range("I8").Value = [double]"5.184741524999481e-13"
$ws.Range("L8").Value = [double]"49.21448745742568"
$ws.Range("M8").Value = "[36.97507977887731, 61.453895135974044]"
$ws.Range("N8").Value = [double]"2.477558158631155e-10"
$ws.Range("O8").Value = [double]"2.477558158631155e-10"
$ws.Range("P8").Value = [double]"-2.603842559732081"
$ws.Range("Q8").Value = "[-2.868000500574466, -2.3396846188896956]"
$ws.Range("T8").Value = [double]"51.58846137251286"
$ws.Range("U8").Value = "[45.21295906447442, 57.9639636805513]"
$ws.Range("X8").Value = [double]"10.70432432432457"
$ws.Range("Y8").Value = [double]"9.618378378378601"
$ws.Range("Z8").Value = [double]"11.79027027027055"
$ws.Range("F9").Value = [double]"25.8300000000006"
$ws.Range("H9").Value = [double]"1.110223024625157e-16"
$ws.Range("I9").Value = [double]"1.110223024625157e-16"
$ws.Range("L9").Value = [double]"49.11188835512958"
$ws.Range("M9").Value = "[41.42809438380347, 56.79568232645569]"
$ws.Range("N9").Value = [double]"2.220446049250313e-16"
$ws.Range("O9").Value = [double]"2.220446049250313e-16"
$ws.Range("P9").Value = [double]"2.899447874484274"
$ws.Range("Q9").Value = "[2.7485004797171966, 3.050395269251351]"
$ws.Range("T9").Value = [double]"51.32207156733372"
$ws.Range("U9").Value = "[46.960911550929566, 55.683231583737864]"
$ws.Range("X9").Value = [double]"13.91045045045077"
$ws.Range("Y9").Value = [double]"13.28990990991022"
$ws.Range("Z9").Value = [double]"14.53099099099133"
$ws.Range("F10").Value = [double]"25.8300000000006"
$ws.Range("H10").Value = [double]"3.752553823233029e-14"
$ws.Range("I10").Value = [double]"3.752553823233029e-14"
$ws.Range("L10").Value = [double]"42.55256864945668"
$ws.Range("M10").Value = "[32.523342050551825, 52.58179524836154]"
$ws.Range("N10").Value = [double]"5.615530263014534e-11"
$ws.Range("O10").Value = [double]"5.615530263014534e-11"
$ws.Range("P10").Value = [double]"2.647868883205812"
$ws.Range("Q10").Value = "[2.408868841491273, 2.8868689249203507]"
$ws.Range("T10").Value = [double]"48.48112299473428"
$ws.Range("U10").Value = "[43.214977937776176, 53.74726805169238]"
$ws.Range("X10").Value = [double]"14.94468468468503"
$ws.Range("Y10").Value = [double]"13.96216216216249"
$ws.Range("Z10").Value = [double]"15.92720720720758"
$ws.Range("F11").Value = [double]"25.8300000000006"
$ws.Range("H11").Value = [double]"1.887379141862766e-15"
$ws.Range("I11").Value = [double]"1.887379141862766e-15"
$ws.Range("L11").Value = [double]"46.67285282575168"
$ws.Range("M11").Value = "[36.23707188380704, 57.10863376769631]"
$ws.Range("N11").Value = [double]"1.235167523816472e-11"
$ws.Range("O11").Value = [double]"1.235167523816472e-11"
$ws.Range("P11").Value = [double]"2.157289850212811"
$ws.Range("Q11").Value = "[1.9308687580621946, 2.383710942363427]"
$ws.Range("T11").Value = [double]"55.20381982947237"
$ws.Range("U11").Value = "[49.728222619109225, 60.67941703983552]"
$ws.Range("X11").Value = [double]"16.96144144144183"
$ws.Range("Y11").Value = [double]"16.030630630631"
$ws.Range("Z11").Value = [double]"17.89225225225267"
$ws.Range("F12").Value = [double]"25.8300000000006"
$ws.Range("H12").Value = [double]"1.110223024625157e-16"
$ws.Range("I12").Value = [double]"1.110223024625157e-16"
$ws.Range("L12").Value = [double]"45.14962648539196"
$ws.Range("M12").Value = "[37.803014455574164, 52.49623851520975]"
$ws.Range("N12").Value = [double]"4.440892098500626e-16"
$ws.Range("O12").Value = [double]"4.440892098500626e-16"
$ws.Range("P12").Value = [double]"1.754763464167271"
$ws.Range("Q12").Value = "[1.5660792207084242, 1.9434477076261185]"
$ws.Range("T12").Value = [double]"47.31768349142437"
$ws.Range("U12").Value = "[42.83057917264393, 51.80478781020481]"
$ws.Range("X12").Value = [double]"18.61621621621665"
$ws.Range("Y12").Value = [double]"17.84054054054095"
$ws.Range("Z12").Value = [double]"19.39189189189235"
$ws.Range("F13").Value = [double]"25.8300000000006"
$ws.Range("H13").Value = [double]"3.552713678800501e-15"
$ws.Range("I13").Value = [double]"3.552713678800501e-15"
$ws.Range("L13").Value = [double]"42.76934903206233"
$ws.Range("M13").Value = "[34.3150898849892, 51.22360817913546]"
$ws.Range("N13").Value = [double]"2.882138971926906e-13"
$ws.Range("O13").Value = [double]"2.882138971926906e-13"
$ws.Range("P13").Value = [double]"1.364816027685655"
$ws.Range("Q13").Value = "[1.1383949355350405, 1.5912371198362694]"
$ws.Range("R13").Value = [double]"8.881784197001252e-16"
$ws.Range("S13").Value = [double]"8.881784197001252e-16"
$ws.Range("T13").Value = [double]"52.84733444018413"
$ws.Range("U13").Value = "[47.7876033085698, 57.907065571798455]"
$ws.Range("X13").Value = [double]"20.21927927927975"
$ws.Range("Y13").Value = [double]"19.28846846846893"
$ws.Range("Z13").Value = [double]"21.15009009009058"
